$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1 from "Round 0" to "Round_0"
$ws.Range("B1").Value = "Round_0"

# Update the selected cell/range to C4
$ws.Range("C4").Select()
